$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.Value = "'51.083.87"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3
$c = $ws.Range("D3")
$c.Value = "'2.960.71"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4
$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$c = $ws.Range("D5")
$c.Value = "'380.15"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "

# Row 6
$c = $ws.Range("D6")
$c.Value = "'102.07"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

# Row 7
$c = $ws.Range("D7")
$c.Value = "'0.544"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.89%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$c = $ws.Range("D9")
$c.Value = "'0.592"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.86%  "

# Row 10
$c = $ws.Range("D10")
$c.Value = "'36.50"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.43%  "

# Row 11
$c = $ws.Range("D11")
$c.Value = "'0.137"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.10%  "

# Row 12
$c = $ws.Range("D12")
$c.Value = "'0.0855"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.32%  "

# Row 13
$c = $ws.Range("D13")
$c.Value = "'3.436.83"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.98%  "

# Row 14
$c = $ws.Range("D14")
$c.Value = "'7.84"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +6.89%  "

# Row 15
$c = $ws.Range("D15")
$c.Value = "'18.32"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.80%  "

# Row 16
$c = $ws.Range("D16")
$c.Value = "'11.80"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +64.85%  "

# Row 17
$c = $ws.Range("D17")
$c.Value = "'2.949.82"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "

# Row 18
$c = $ws.Range("D18")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.44%  "

# Row 19
$c = $ws.Range("D19")
$c.Value = "'51.182.39"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20
$c = $ws.Range("D20")
$c.Value = "'3.11"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.53%  "

# Row 21
$c = $ws.Range("D21")
$c.Value = "'12.41"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.21%  "

# Row 22
$c = $ws.Range("D22")
$c.Value = "'0.0₃0960"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23
$c = $ws.Range("D23")
$c.Value = "'70.04"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.54%  "

# Row 24
$c = $ws.Range("D24")
$c.Value = "'3.29"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +14.35%  "

# Row 25
$c = $ws.Range("D25")
$c.Value = "'266.94"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "

# Row 26
$c = $ws.Range("D26")
$c.Value = "'7.92"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -4.66%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D27")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D28")
$c.Value = "'7.11"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -10.32%  "

# Row 29
$c = $ws.Range("D29")
$c.Value = "'0.166"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "

# Row 30
$c = $ws.Range("D30")
$c.Value = "'25.87"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.17%  "

# Row 31
$ws.Range("E31").Value = "  -1.93%  "

# Row 32
$c = $ws.Range("D32")
$c.Value = "'10.29"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.48%  "

# Row 33
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D33")
$c.Value = "'51.25"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.92%  "

# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D34")
$c.Value = "'34.37"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.62%  "

# Row 35
$c = $ws.Range("D35")
$c.Value = "'2.06"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.89%  "

# Row 36
$c = $ws.Range("D36")
$c.Value = "'0.0435"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.56%  "

# Row 37
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$c = $ws.Range("D38")
$c.Value = "'3.24"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +8.78%  "

# Row 39
$ws.Range("E39").Value = "  +1.22%  "

# Row 40
$c = $ws.Range("D40")
$c.Value = "'1.83"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.76%  "

# Row 41
$c = $ws.Range("D41")
$c.Value = "'16.49"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D42")
$c.Value = "'2.51"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.97%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D43")
$c.Value = "'124.76"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +3.56%  "

# Row 44
$c = $ws.Range("D44")
$c.Value = "'21.52"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.76%  "

# Row 45
$c = $ws.Range("D45")
$c.Value = "'3.54"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +9.30%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D46")
$c.Value = "'0.273"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -5.45%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Range("D47")
$c.Value = "'2.38"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.10%  "

# Row 48
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D48")
$c.Value = "'2.02"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "

# Row 49
$c = $ws.Range("D49")
$c.Value = "'2.052.51"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.02%  "

# Row 50
$c = $ws.Range("D50")
$c.Value = "'0.0320"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -6.16%  "

# Row 51
$c = $ws.Range("D51")
$c.Value = "'5.42"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +7.52%  "
